# Updated league table for GW25.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 24 (row 22) results just came in - fill in Eren/Mert/Arda's scores.
$ws.Range("B22").Value = 92
$ws.Range("C22").Value = 93
$ws.Range("D22").Value = 98

# Carry over the same (blank/no-op) cell formatting used by the rows above
# it in the table (rows 20-21) so the new row matches its neighbours.
$ws.Range("B22:D22").Style = $ws.Range("B21:D21").Style

# Leave the cursor where the author last clicked while reviewing the update.
$ws.Range("I21").Select()
